# Insert a new weekly data row at row 568 (pushing existing rows 568..603 down to 569..604)
# and populate it with the new "Inferno" Ají record for Limache.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("568:568").Insert()

$ws.Range("A568").Value = 3
$ws.Range("B568").Value = "Femacal de La Calera"
$ws.Range("C568").Value = "Coquimbo"
$ws.Range("D568").Value = (Get-Date -Year 2023 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E568").Value = 5
$ws.Range("F568").Value = 100112021
$ws.Range("G568").Value = "Ají"
$ws.Range("H568").Value = "Inferno"
$ws.Range("I568").Value = "Primera"
$ws.Range("J568").Value = 76
$ws.Range("K568").Value = 19000
$ws.Range("L568").Value = 20000
$ws.Range("M568").Value = 19500
$ws.Range("N568").Value = "$/caja 15 kilos"
$ws.Range("O568").Value = "Limache"
$ws.Range("P568").Value = 1300
$ws.Range("Q568").Value = 15
$ws.Range("R568").Value = "Hortaliza"
